# Weekly price update: insert a new price record for week of 2022-07-04
# (row 100) into the "Agrícola del Norte S.A. de Arica - Palta" sheet,
# pushing the existing rows 100-105 down to 101-106.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 100 (shifts rows 100:105 -> 101:106)
$ws.Rows.Item(100).Insert()

# Populate the newly inserted row with this week's Hass / Primera data point
$ws.Range("A100").Value = 1
$ws.Range("B100").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C100").Value = "Arica y Parinacota"
$ws.Range("D100").Value = 44746
$ws.Range("E100").Value = 15
$ws.Range("F100").Value = "Fruta"
$ws.Range("G100").Value = 100106
$ws.Range("H100").Value = "Oleaginosos"
$ws.Range("I100").Value = 100106002
$ws.Range("J100").Value = "Palta"
$ws.Range("K100").Value = "Hass"
$ws.Range("L100").Value = "Primera"
$ws.Range("M100").Value = 400
$ws.Range("N100").Value = 14000
$ws.Range("O100").Value = 15000
$ws.Range("P100").Value = 14500
$ws.Range("Q100").Value = "`$/bandeja 10 kilos"
$ws.Range("R100").Value = "Perú"
$ws.Range("S100").Value = 1450
$ws.Range("T100").Value = 10
